# Adapt column header formatting to respective input file names:
#   "<name>_old"  -> "<name>_FV2210"
#   "<name>_new"  -> "<name>_FV2304"
# Also wrap the data range in an Excel Table ("Table1") with autofilter,
# and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename the header row (row 1, columns A:U) -------------------------
$headers = @(
  "Segmentname_FV2210", "Segmentgruppe_FV2210", "Segment_FV2210", "Datenelement_FV2210",
  "Segment ID_FV2210", "Code_FV2210", "Qualifier_FV2210", "Beschreibung_FV2210",
  "Bedingungsausdruck_FV2210", "Bedingung_FV2210",
  "diff",
  "Segmentname_FV2304", "Segmentgruppe_FV2304", "Segment_FV2304", "Datenelement_FV2304",
  "Segment ID_FV2304", "Code_FV2304", "Qualifier_FV2304", "Beschreibung_FV2304",
  "Bedingungsausdruck_FV2304", "Bedingung_FV2304"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# --- 2. Turn the used range into a real Excel Table (Table1) ---------------
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U70"), $null, 1)
$lo.Name = "Table1"
$lo.TableStyle = ""

# --- 3. Freeze the header row (split after row 1) ---------------------------
# Selecting A2 first (the cell just below the split) before freezing makes
# Excel record pane state="frozen" (as opposed to "frozenSplit"); the
# selection is then moved back to A1 to match the usual resting state.
$ws.Activate() | Out-Null
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A1").Select() | Out-Null
